{"js": "// Append the \"WHY / HOW / WHAT\" section to the landing page content.\n//\n// The new content goes right after the last paragraph of body text\n// (\"...you will receive the funding you deserve!\") and right before the\n// pre-existing trailing empty paragraph, producing this shape:\n//\n//   <empty paragraph>\n//   WHY- ... you are.\n//   <empty paragraph>\n//   HOW- ... look for funding.\n//   WHAT-We're a platform ... good enough for us.\n//   <original trailing empty paragraph>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the paragraph that currently ends the body copy.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"you will receive the funding you deserve!\") !== -1) {\n    targetIndex = i;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the 'you deserve!' paragraph to anchor the new content.\");\n}\n\nconst anchorParagraph = paragraphs.items[targetIndex];\n\nconst whyText =\n  \"WHY- We believe that an idea should never have to die because of a lack of money. \" +\n  \"We believe that passion should never have to die because of a lack of money. We understand \" +\n  \"that approaching investors is not only difficult, but a rather tedious process. And since, \" +\n  \"we\\u2019ve brought up the word \\u201cinvestor\\u201d, what really makes a person an investor? \" +\n  \"A million dollars in the bank? A swanky car and a fancy number plate? We don\\u2019t think so. \" +\n  \"If you\\u2019re willing to invest $3000, you\\u2019re still an investor, aren\\u2019t you? \" +\n  \"In our eyes, we definitely believe you are.\";\n\nconst howText =\n  \"HOW- If there\\u2019s anything that we\\u2019ve learnt from the COVID-19 quarantine, it is the \" +\n  \"importance of technology and the role it plays in keeping us connected. A simple, yet efficient \" +\n  \"user interface to not only make sure your passion will be funded, but also for you to have a \" +\n  \"seamless experience while you look for funding.\";\n\nconst whatText =\n  \"WHAT-We\\u2019re a platform that connects entrepreneurs to people who want to invest. As an \" +\n  \"entrepreneur(ship), you don\\u2019t need to be the next Google. You might want to start a cake \" +\n  \"shop in your locality because the current one is no good. As an investor, you don\\u2019t need \" +\n  \"to be the richest man on the planet. You might want to invest amounts near Rs. 1 lac, and \" +\n  \"that\\u2019s good enough for us.\";\n\n// Each insertParagraph(..., Word.InsertLocation.after) lands immediately\n// after the anchor paragraph (before the old trailing empty paragraph), so\n// inserting in reverse order yields the correct final reading order:\n//   anchor, <empty>, WHY, <empty>, HOW, WHAT, <old trailing empty>\nanchorParagraph.insertParagraph(whatText, Word.InsertLocation.after);\nanchorParagraph.insertParagraph(howText, Word.InsertLocation.after);\nanchorParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nanchorParagraph.insertParagraph(whyText, Word.InsertLocation.after);\nanchorParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Append the \"WHY / HOW / WHAT\" section to the landing page content.\n#\n# The new content goes right after the last paragraph of body text\n# (\"...you will receive the funding you deserve!\") and right before the\n# pre-existing trailing empty paragraph, following this shape:\n#\n#   <empty paragraph>\n#   WHY- ... you are.\n#   <empty paragraph>\n#   HOW- ... look for funding.\n#   WHAT-We're a platform ... good enough for us.\n#   <original trailing empty paragraph>\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that currently ends the body copy.\n$paragraphs = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $t = $paragraphs.Item($i).Range.Text\n    if ($t -like '*you will receive the funding you deserve!*') {\n        $targetIndex = $i\n    }\n}\n\n$lastPara = $paragraphs.Item($targetIndex)\n$lastRange = $lastPara.Range\n\n# Insert one fresh paragraph right after it; we'll fill this paragraph\n# (and the ones that follow) with the whole new block in a single\n# Range.Text assignment so every CR becomes its own paragraph, inserted\n# immediately before the document's original trailing empty paragraph.\n$lastRange.InsertParagraphAfter()\n$insertedRange = $d.Paragraphs.Item($targetIndex + 1).Range\n\n$cr = [char]13\n\n$whyText = 'WHY- We believe that an idea should never have to die because of a lack of money. We believe that passion should never have to die because of a lack of money. We understand that approaching investors is not only difficult, but a rather tedious process. And since, we\u2019ve brought up the word \u201cinvestor\u201d, what really makes a person an investor? A million dollars in the bank? A swanky car and a fancy number plate? We don\u2019t think so. If you\u2019re willing to invest $3000, you\u2019re still an investor, aren\u2019t you? In our eyes, we definitely believe you are.'\n\n$howText = 'HOW- If there\u2019s anything that we\u2019ve learnt from the COVID-19 quarantine, it is the importance of technology and the role it plays in keeping us connected. A simple, yet efficient user interface to not only make sure your passion will be funded, but also for you to have a seamless experience while you look for funding.'\n\n$whatText = 'WHAT-We\u2019re a platform that connects entrepreneurs to people who want to invest. As an entrepreneur(ship), you don\u2019t need to be the next Google. You might want to start a cake shop in your locality because the current one is no good. As an investor, you don\u2019t need to be the richest man on the planet. You might want to invest amounts near Rs. 1 lac, and that\u2019s good enough for us.'\n\n$parts = @('', $whyText, '', $howText, $whatText)\n$fullBlock = [string]::Join($cr, $parts)\n\n$insertedRange.Text = $fullBlock\n\nWrite-Output \"Inserted WHY/HOW/WHAT block after paragraph $targetIndex.\"\n"}
